$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest scraped values for the "Price" (D) and "Volume(1h)" (E) columns of the
# cryptocurrency table, refreshed by the scheduled GitHub Actions scraper run.
# NumberFormat is forced to Text before each write (and the style reset to
# "Normal" afterwards) so that numeric-looking values such as "603.03" are
# stored verbatim as text -- exactly as scraped -- instead of being coerced by
# Excel into floating point numbers (which would corrupt values like "7.80" by
# dropping the trailing zero) and without leaving a residual cell style behind.
$updates = @{
    "D2" = "66.214.46"
    "E2" = "  -0.13%  "
    "D3" = "3.547.79"
    "E3" = "  -0.17%  "
    "D5" = "603.03"
    "E5" = "  -0.35%  "
    "D6" = "146.39"
    "E6" = "  +1.54%  "
    "D7" = "3.548.89"
    "E7" = "  -0.13%  "
    "E8" = "  -0.23%  "
    "D9" = "0.497"
    "E9" = "  +1.50%  "
    "E10" = "  -1.81%  "
    "D11" = "7.80"
    "E11" = "  -0.42%  "
    "E12" = "  -0.97%  "
    "D13" = "4.144.51"
    "E13" = "  -0.34%  "
    "D14" = "0.0000202"
    "E14" = "  -2.31%  "
    "D15" = "29.05"
    "E15" = "  -3.41%  "
    "D16" = "3.543.93"
    "E16" = "  -0.43%  "
    "E17" = "  +1.69%  "
    "D18" = "66.165.46"
    "E18" = "  -0.38%  "
    "D19" = "11.07"
    "E19" = "  -3.33%  "
    "D20" = "6.23"
    "E20" = "  +0.82%  "
    "D21" = "14.67"
    "D22" = "418.86"
    "E22" = "  -2.75%  "
    "D23" = "0.602"
    "E23" = "  -1.13%  "
    "D24" = "77.78"
    "E24" = "  -2.25%  "
    "D25" = "3.682.53"
    "E25" = "  -0.44%  "
    "E26" = "  +0.08%  "
    "E27" = "  -2.30%  "
    "E28" = "  -0.30%  "
    "E29" = "  -1.46%  "
    "D30" = "7.83"
    "E30" = "  -1.54%  "
    "E31" = "  +0.06%  "
    "D32" = "3.541.21"
    "E32" = "  -0.24%  "
    "E33" = "  +1.75%  "
    "D34" = "24.43"
    "D36" = "7.58"
    "E36" = "  -2.94%  "
    "D37" = "1.31"
    "E37" = "  -9.45%  "
    "D38" = "174.22"
    "E38" = "  -1.10%  "
    "E39" = "  -7.12%  "
    "D40" = "5.27"
    "E40" = "  -5.60%  "
    "D41" = "0.0824"
    "E41" = "  -2.65%  "
    "D42" = "5.10"
    "E42" = "  -1.66%  "
    "D43" = "0.862"
    "E43" = "  -2.81%  "
    "D44" = "45.69"
    "E44" = "  -0.60%  "
    "E45" = "  -5.77%  "
    "E46" = "  -0.03%  "
    "D47" = "2.43"
    "E47" = "  -3.81%  "
    "D48" = "7.13"
    "E48" = "  +0.08%  "
    "D49" = "22.87"
    "E49" = "  -1.55%  "
    "D50" = "1.11"
    "E50" = "  -7.23%  "
    "D51" = "23.34"
    "E51" = "  -6.97%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
